$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Задание рассмотрено на заседании цикловой комиссии по УГС 09.00.00
#     Информатика и вычислительная техника" (highlighted yellow)
#    -> "...по специальности 09.02.07 Информационные технологии и
#     программирование" (no highlight)
# ---------------------------------------------------------------------------
$p47 = $d.Paragraphs(47)
$p47start = $p47.Range.Start
$oldTail = $d.Range($p47start + 54, $p47start + 103)
$oldTail.Text = ""
$insPoint = $d.Range($p47start + 54, $p47start + 54)
$insPoint.InsertAfter("специальности 09.02.07 Информационные технологии и программирование")

# ---------------------------------------------------------------------------
# 2) "Протокол № [tab]от «[tab]»[tab]2022 г." ->
#    "Протокол № 6 от «21» апреля 2022 г."
# ---------------------------------------------------------------------------
$p50 = $d.Paragraphs(50)
$p50start = $p50.Range.Start

# third tab (offset 18) -> " апреля "
$t3 = $d.Range($p50start + 18, $p50start + 19)
$t3.Text = ""
$ip3 = $d.Range($p50start + 18, $p50start + 18)
$ip3.InsertAfter(" апреля ")

# second tab (offset 16) -> "21"
$t2 = $d.Range($p50start + 16, $p50start + 17)
$t2.Text = ""
$ip2 = $d.Range($p50start + 16, $p50start + 16)
$ip2.InsertAfter("21")

# first tab (offset 11) -> "6 "
$t1 = $d.Range($p50start + 11, $p50start + 12)
$t1.Text = ""
$ip1 = $d.Range($p50start + 11, $p50start + 11)
$ip1.InsertAfter("6 ")

# ---------------------------------------------------------------------------
# 3) "Председатель цикловой комиссии" -> "Председатель цикловой комиссии "
# ---------------------------------------------------------------------------
$p51 = $d.Paragraphs(51)
$p51start = $p51.Range.Start
$ipPred = $d.Range($p51start + 30, $p51start + 30)
$ipPred.InsertAfter(" ")

# ---------------------------------------------------------------------------
# 4) insert a space run (same formatting as the tabs) right before "О.О. "
#    (after the last of the 5 underline tabs). Paragraph 51 grew by 1 char
#    because of edit (3) above, so the 5 tabs now sit at offsets 31..35.
# ---------------------------------------------------------------------------
$p51b = $d.Paragraphs(51)
$p51bstart = $p51b.Range.Start
$ipSpace = $d.Range($p51bstart + 36, $p51bstart + 36)
$ipSpace.InsertAfter(" ")
$spaceRange = $d.Range($p51bstart + 36, $p51bstart + 37)
$spaceRange.Font.Underline = 1

# ---------------------------------------------------------------------------
# 5) "руководитель курсового проектирования" (р + уководитель...) ->
#    "Руководитель курсового проектирования " (Р + уководитель...<space>)
# ---------------------------------------------------------------------------
$p52 = $d.Paragraphs(52)
$p52start = $p52.Range.Start
$firstLetter = $d.Range($p52start, $p52start + 1)
$firstLetter.Text = "Р"
$p52b = $d.Paragraphs(52)
$p52bstart = $p52b.Range.Start
$ipRuk = $d.Range($p52bstart + 38, $p52bstart + 38)
$ipRuk.InsertAfter(" ")

# ---------------------------------------------------------------------------
# 6) "«___» ________ 2022 г." -> "«17» марта 2022 г."
# ---------------------------------------------------------------------------
$p56 = $d.Paragraphs(56)
$p56start = $p56.Range.Start
$oldDate = $d.Range($p56start, $p56start + 22)
$oldDate.Text = ""
$ipDate = $d.Range($p56start, $p56start)
$ipDate.InsertAfter("«17» марта 2022 г.")

# ---------------------------------------------------------------------------
# 7) add a new, empty paragraph at the very end of the body (after the
#    "«23» мая 2022 г." paragraph), with the same pPr/tabs.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endOfBody = $lastPara.Range.End
$endRange = $d.Range($endOfBody, $endOfBody)
$endRange.InsertParagraphAfter()
